# Data retrieved - Sat May 15 18:37:42 UTC 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17, col A (date) gets a tiny float update (re-fetched timestamp).
$ws.Cells.Item(17, 1).Value = 44330.77832976389

# New row 18 with the freshly retrieved job numbers.
$ws.Cells.Item(18, 1).Value = 44331.77618375275
$ws.Cells.Item(18, 2).Value = 74352
$ws.Cells.Item(18, 3).Value = 62593
$ws.Cells.Item(18, 4).Value = 3222
$ws.Cells.Item(18, 5).Value = 2110
$ws.Cells.Item(18, 6).Value = 1492
$ws.Cells.Item(18, 7).Value = 19331
$ws.Cells.Item(18, 8).Value = 1303
$ws.Cells.Item(18, 9).Value = 868
$ws.Cells.Item(18, 10).Value = 204

# Column A uses a custom date/time number format (style index 2 in
# before.xlsx) - apply the same number format to the new row's date cell,
# matching the rest of column A.
$ws.Cells.Item(18, 1).NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
